$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The rights to both movies have been successfully acquired.`n"
$ws.Range("D2").Value = "both_movies, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday.`n"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision to acquire the rights for ""Barbie"" has been successfully recorded.`n"
$ws.Range("D4").Value = "Barbie_was_selected, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision indicates that no agreement was reached regarding the movie selection for Friday.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" as the movie for Friday's assembly.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision was made to not select a movie for Friday's showing.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The final decision has been recorded as no decision about the movie for Friday's showing.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision remains that no movie was selected.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision to acquire the rights to both movies has been successfully recorded.`n"
$ws.Range("D10").Value = "both_movies, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The conversation did not result in a decision regarding the movie to be shown on Friday, so no acquisition of movie rights will take place.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision to acquire the rights for ""Barbie"" has been recorded.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision process concluded without a definitive choice for Friday's movie.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for ""Barbie"" have been acquired for showing on Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision to select ""Barbie"" as the movie for Friday has been confirmed.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision resulted in no selection for the movie to be shown on Friday.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded successfully. The movie ""Barbie"" will be shown on Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not reached, resulting in no selection being made.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded as ""no decision"" regarding the movie for Friday.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as no selection for Friday's movie.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The committee did not reach a decision about which movie to show on Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been registered with no movie selected for Friday.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision has been recorded, and no selection was made for Friday's movie.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded successfully: ""Barbie"" will be the movie acquired for Friday's showing.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The rights to both movies will be acquired.`n"
$ws.Range("D28").Value = "both_movies, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision has been recorded and the rights for ""Barbie"" have been acquired.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision process regarding the movie to show on Friday resulted in no decision being made.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been recorded and no movie has been selected for Friday.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" to be shown on Friday.`n"
$ws.Range("D33").Value = "Barbie_was_selected, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision process ended without reaching a conclusion on which movie to show on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement, and thus no movie has been selected.`n"
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been recorded as no definitive choice for the movie on Friday.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: No decision has been made regarding the movie to be shown on Friday.`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement on which movie to show on Friday.`n"
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for ""Barbie.""`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision to acquire the rights for ""Oppenheimer"" has been successfully recorded.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision-making process concluded without a definitive choice for Friday's movie, resulting in no decision being made.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for ""Barbie.""`n"
$ws.Range("D43").Value = "Barbie_was_selected, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no final agreement was reached on the movie selection for Friday.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" as the movie to be shown on Friday.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached on which movie to show on Friday.`n"
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie.""`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision process did not lead to an agreement on which movie to show, resulting in no decision being made.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection has been made for Friday.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision to show a movie on Friday could not be made during the meeting.`n"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie.""`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded to select ""Barbie"" as the movie to be shown on Friday.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The conclusion of the meeting was that no decision was made regarding the movie for Friday.`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision regarding which movie to acquire for Friday has resulted in no agreement.`n"
